# ---------------------------------------------------------------
# Update existing odds values in rows 2-10 (per the commit diff)
# ---------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 2.44
$ws.Range("G2").Value = 2.48
$ws.Range("I2").Value = 3.6
$ws.Range("J2").Value = 3.2
$ws.Range("K2").Value = 3.25
$ws.Range("T2").Value = 1.93
$ws.Range("X2").Value = 9.8
$ws.Range("Y2").Value = 11.5
$ws.Range("Z2").Value = 23
$ws.Range("AB2").Value = 8.8
$ws.Range("AC2").Value = 7
$ws.Range("AO2").Value = 55

# Row 3
$ws.Range("G3").Value = 1.79
$ws.Range("H3").Value = 6
$ws.Range("I3").Value = 7
$ws.Range("J3").Value = 3.6
$ws.Range("K3").Value = 3.85
$ws.Range("L3").Value = 1.56
$ws.Range("N3").Value = 2.92
$ws.Range("S3").Value = 5
$ws.Range("T3").Value = 2.12
$ws.Range("U3").Value = 1.7
$ws.Range("W3").Value = 2.06
$ws.Range("Y3").Value = 1000
$ws.Range("Z3").Value = 1000
$ws.Range("AD3").Value = 1000
$ws.Range("AF3").Value = 9.4
$ws.Range("AG3").Value = 10.5
$ws.Range("AK3").Value = 1000

# Row 4
$ws.Range("I4").Value = 2.88
$ws.Range("N4").Value = 2.94
$ws.Range("P4").Value = 1.67
$ws.Range("Q4").Value = 2.3
$ws.Range("V4").Value = 1.53
$ws.Range("AC4").Value = 8.4
$ws.Range("AD4").Value = 14.5
$ws.Range("AO4").Value = 40

# Row 5
$ws.Range("T5").Value = 2.14
$ws.Range("AB5").Value = 8

# Row 6
$ws.Range("F6").Value = 1.8
$ws.Range("G6").Value = 1.9
$ws.Range("I6").Value = 6
$ws.Range("P6").Value = 1.72
$ws.Range("Q6").Value = 2.12
$ws.Range("S6").Value = 4
$ws.Range("T6").Value = 1.98
$ws.Range("V6").Value = 1.2
$ws.Range("W6").Value = 2.1
$ws.Range("X6").Value = 11.5
$ws.Range("AA6").Value = 160
$ws.Range("AB6").Value = 7.8
$ws.Range("AE6").Value = 90
$ws.Range("AM6").Value = 170
$ws.Range("AN6").Value = 16
$ws.Range("AO6").Value = 130

# Row 7
$ws.Range("G7").Value = 1.1
$ws.Range("H7").Value = 30
$ws.Range("J7").Value = 13
$ws.Range("L7").Value = 1.2
$ws.Range("N7").Value = 8.4
$ws.Range("O7").Value = 1.1
$ws.Range("P7").Value = 3.45
$ws.Range("Q7").Value = 1.32
$ws.Range("R7").Value = 2.04
$ws.Range("S7").Value = 1.81
$ws.Range("T7").Value = 3.4
$ws.Range("U7").Value = 1.4
$ws.Range("W7").Value = 11
$ws.Range("Y7").Value = 160
$ws.Range("AB7").Value = 15
$ws.Range("AN7").Value = 2.36

# Row 8
$ws.Range("F8").Value = 1.65
$ws.Range("G8").Value = 1.72
$ws.Range("J8").Value = 3.6
$ws.Range("K8").Value = 3.85
$ws.Range("M8").Value = 1.09
$ws.Range("P8").Value = 1.63
$ws.Range("S8").Value = 4.7
$ws.Range("U8").Value = 1.68
$ws.Range("W8").Value = 2.38

# Row 9
$ws.Range("Q9").Value = 2.06
$ws.Range("U9").Value = 1.71
$ws.Range("V9").Value = 1.09

# Row 10
$ws.Range("F10").Value = 1.33
$ws.Range("G10").Value = 1.38
$ws.Range("H10").Value = 12
$ws.Range("I10").Value = 15.5
$ws.Range("J10").Value = 5.1
$ws.Range("K10").Value = 6
$ws.Range("N10").Value = 3.95
$ws.Range("P10").Value = 2.02
$ws.Range("Q10").Value = 1.78
$ws.Range("R10").Value = 1.4
$ws.Range("S10").Value = 3
$ws.Range("T10").Value = 2.26
$ws.Range("U10").Value = 1.65
$ws.Range("V10").Value = 1.06
$ws.Range("W10").Value = 3.6
$ws.Range("X10").Value = 22
$ws.Range("Z10").Value = 160
$ws.Range("AA10").Value = 1000
$ws.Range("AB10").Value = 9.2
$ws.Range("AC10").Value = 15.5
$ws.Range("AD10").Value = 1000
$ws.Range("AE10").Value = 330
$ws.Range("AF10").Value = 7.8
$ws.Range("AG10").Value = 11.5
$ws.Range("AH10").Value = 42
$ws.Range("AI10").Value = 250
$ws.Range("AJ10").Value = 10.5
$ws.Range("AK10").Value = 16.5
$ws.Range("AM10").Value = 290
$ws.Range("AN10").Value = 6.4

# ---------------------------------------------------------------
# Insert a new row at position 11 (shifts old row 11 -> row 12,
# and extends the used range to A1:AO12)
# ---------------------------------------------------------------
$ws.Rows(11).Insert()

# ---------------------------------------------------------------
# Populate the brand-new row 11: Honduras Liga Nacional match
# ---------------------------------------------------------------
$ws.Range("A11").Value = "Honduras Liga Nacional"
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "2025-10-14"
$ws.Range("B11").ClearFormats()
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "22:30:00"
$ws.Range("C11").ClearFormats()
$ws.Range("D11").Value = "Atletico Choloma"
$ws.Range("E11").Value = "Lobos UPNFM"
$ws.Range("F11").Value = 1.04
$ws.Range("G11").Value = 1000
$ws.Range("H11").Value = 1.04
$ws.Range("I11").Value = 1000
$ws.Range("J11").Value = 1.02
$ws.Range("K11").Value = 950
$ws.Range("L11").Value = 1.01
$ws.Range("M11").Value = 1.01
$ws.Range("N11").Value = 1.3
$ws.Range("O11").Value = 1.22
$ws.Range("P11").Value = 1.3
$ws.Range("Q11").Value = 1.22
$ws.Range("R11").Value = 1.18
$ws.Range("S11").Value = 1.22
$ws.Range("T11").Value = 1.01
$ws.Range("U11").Value = 1.01
$ws.Range("V11").Value = 1.01
$ws.Range("W11").Value = 1.01
$ws.Range("X11").Value = 1000
$ws.Range("Y11").Value = 1000
$ws.Range("Z11").Value = 1000
$ws.Range("AA11").Value = 1000
$ws.Range("AB11").Value = 1000
$ws.Range("AC11").Value = 1000
$ws.Range("AD11").Value = 1000
$ws.Range("AE11").Value = 1000
$ws.Range("AF11").Value = 1000
$ws.Range("AG11").Value = 1000
$ws.Range("AH11").Value = 1000
$ws.Range("AI11").Value = 1000
$ws.Range("AJ11").Value = 1000
$ws.Range("AK11").Value = 1000
$ws.Range("AL11").Value = 1000
$ws.Range("AM11").Value = 1000
$ws.Range("AN11").Value = 1000
$ws.Range("AO11").Value = 1000
# ---------------------------------------------------------------
# Row 12 (shifted from old row 11): El Salvador vs Guatemala --
# text columns A-E already carried over correctly by the Insert;
# only the odds columns F:AO need new values
# ---------------------------------------------------------------
$ws.Range("F12").Value = 2.58
$ws.Range("G12").Value = 2.88
$ws.Range("H12").Value = 2.84
$ws.Range("I12").Value = 3.2
$ws.Range("J12").Value = 3.1
$ws.Range("K12").Value = 3.55
$ws.Range("L12").Value = 1.52
$ws.Range("M12").Value = 1.1
$ws.Range("N12").Value = 2.92
$ws.Range("O12").Value = 1.44
$ws.Range("P12").Value = 1.64
$ws.Range("Q12").Value = 2.26
$ws.Range("R12").Value = 1.24
$ws.Range("S12").Value = 4.7
$ws.Range("T12").Value = 1.9
$ws.Range("U12").Value = 1.87
$ws.Range("V12").Value = 1.45
$ws.Range("W12").Value = 1.53
$ws.Range("X12").Value = 1000
$ws.Range("Y12").Value = 1000
$ws.Range("Z12").Value = 24
$ws.Range("AA12").Value = 1000
$ws.Range("AB12").Value = 1000
$ws.Range("AC12").Value = 7.8
$ws.Range("AD12").Value = 15
$ws.Range("AE12").Value = 1000
$ws.Range("AF12").Value = 1000
$ws.Range("AG12").Value = 1000
$ws.Range("AH12").Value = 1000
$ws.Range("AI12").Value = 1000
$ws.Range("AJ12").Value = 1000
$ws.Range("AK12").Value = 1000
$ws.Range("AL12").Value = 1000
$ws.Range("AM12").Value = 180
$ws.Range("AN12").Value = 1000
$ws.Range("AO12").Value = 1000